$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1954022988505747
$ws.Range("C2").Value = 0.5545977011494253
$ws.Range("J2").Value = 0.02011494252873563
$ws.Range("P2").Value = 0.1494252873563219
$ws.Range("S2").Value = 0.08045977011494253

# Row 3
$ws.Range("B3").Value = 0.004629629629629629
$ws.Range("C3").Value = 0.03703703703703703
$ws.Range("J3").Value = 0.04166666666666666
$ws.Range("P3").Value = 0.7592592592592593
$ws.Range("S3").Value = 0.1574074074074074

# Row 4
$ws.Range("J4").Value = 0.08333333333333333
$ws.Range("O4").Value = 0.01666666666666667
$ws.Range("P4").Value = 0.5333333333333333
$ws.Range("S4").Value = 0.3666666666666666

# Row 6
$ws.Range("B6").Value = 0.07773851590106007
$ws.Range("D6").Value = 0.01413427561837456
$ws.Range("F6").Value = 0.04593639575971731
$ws.Range("J6").Value = 0.2049469964664311
$ws.Range("O6").Value = 0.02826855123674912
$ws.Range("Q6").Value = 0.1908127208480565
$ws.Range("R6").Value = 0.07773851590106007
$ws.Range("S6").Value = 0.3604240282685512

# Row 7
$ws.Range("B7").Value = 0.1382978723404255
$ws.Range("D7").Value = 0.01063829787234043
$ws.Range("F7").Value = 0.05319148936170213
$ws.Range("J7").Value = 0.1117021276595745
$ws.Range("O7").Value = 0.02127659574468085
$ws.Range("Q7").Value = 0.1914893617021277
$ws.Range("R7").Value = 0.04787234042553191
$ws.Range("S7").Value = 0.425531914893617

# Row 8
$ws.Range("B8").Value = 0.07006369426751592
$ws.Range("D8").Value = 0.0148619957537155
$ws.Range("E8").Value = 0.002123142250530786
$ws.Range("F8").Value = 0.07218683651804671
$ws.Range("J8").Value = 0.1146496815286624
$ws.Range("O8").Value = 0.02123142250530785
$ws.Range("Q8").Value = 0.2038216560509554
$ws.Range("R8").Value = 0.07006369426751592
$ws.Range("S8").Value = 0.4309978768577495

# Row 9
$ws.Range("B9").Value = 0.1235521235521236
$ws.Range("D9").Value = 0.02316602316602316
$ws.Range("F9").Value = 0.08108108108108109
$ws.Range("J9").Value = 0.1042471042471042
$ws.Range("O9").Value = 0.01544401544401544
$ws.Range("Q9").Value = 0.2123552123552123
$ws.Range("R9").Value = 0.07722007722007722
$ws.Range("S9").Value = 0.3629343629343629

# Row 10
$ws.Range("B10").Value = 0.1137473831123517
$ws.Range("D10").Value = 0.02651779483600837
$ws.Range("E10").Value = 0.001395673412421493
$ws.Range("F10").Value = 0.08792742498255408
$ws.Range("J10").Value = 0.1053733426378227
$ws.Range("O10").Value = 0.02512212142358688
$ws.Range("Q10").Value = 0.2135380321004885
$ws.Range("R10").Value = 0.07466852756454989
$ws.Range("S10").Value = 0.3517096999302163

# Row 11
$ws.Range("G11").Value = 0.1206896551724138
$ws.Range("J11").Value = 0.1206896551724138
$ws.Range("K11").Value = 0.1551724137931035
$ws.Range("L11").Value = 0.5827586206896552
$ws.Range("S11").Value = 0.02068965517241379

# Row 12
$ws.Range("G12").Value = 0.7472527472527473
$ws.Range("J12").Value = 0.1483516483516484
$ws.Range("K12").Value = 0.01098901098901099
$ws.Range("L12").Value = 0.06043956043956044
$ws.Range("S12").Value = 0.03296703296703297

# Row 13
$ws.Range("G13").Value = 0.5348837209302325
$ws.Range("J13").Value = 0.3720930232558139
$ws.Range("S13").Value = 0.09302325581395349

# Row 15
$ws.Range("F15").Value = 0.0132890365448505
$ws.Range("H15").Value = 0.1661129568106312
$ws.Range("I15").Value = 0.07308970099667775
$ws.Range("J15").Value = 0.3654485049833887
$ws.Range("K15").Value = 0.04983388704318937
$ws.Range("M15").Value = 0.0132890365448505
$ws.Range("O15").Value = 0.08637873754152824
$ws.Range("S15").Value = 0.2325581395348837

# Row 16
$ws.Range("F16").Value = 0.02654867256637168
$ws.Range("H16").Value = 0.1327433628318584
$ws.Range("I16").Value = 0.1061946902654867
$ws.Range("J16").Value = 0.415929203539823
$ws.Range("K16").Value = 0.1194690265486726
$ws.Range("M16").Value = 0.02654867256637168
$ws.Range("O16").Value = 0.04424778761061947
$ws.Range("S16").Value = 0.1283185840707965

# Row 17
$ws.Range("F17").Value = 0.02154398563734291
$ws.Range("H17").Value = 0.1813285457809695
$ws.Range("I17").Value = 0.1202872531418312
$ws.Range("J17").Value = 0.4434470377019749
$ws.Range("K17").Value = 0.0718132854578097
$ws.Range("M17").Value = 0.01795332136445242
$ws.Range("O17").Value = 0.06822262118491922
$ws.Range("S17").Value = 0.07540394973070018

# Row 18
$ws.Range("F18").Value = 0.01538461538461539
$ws.Range("H18").Value = 0.1743589743589744
$ws.Range("I18").Value = 0.1230769230769231
$ws.Range("J18").Value = 0.4461538461538462
$ws.Range("K18").Value = 0.08717948717948718
$ws.Range("M18").Value = 0.005128205128205128
$ws.Range("O18").Value = 0.09743589743589744
$ws.Range("S18").Value = 0.05128205128205128

# Row 19
$ws.Range("F19").Value = 0.01675552170601675
$ws.Range("H19").Value = 0.1964965727341965
$ws.Range("I19").Value = 0.09520182787509521
$ws.Range("J19").Value = 0.38994668697639
$ws.Range("K19").Value = 0.1096725057121097
$ws.Range("M19").Value = 0.01827875095201828
$ws.Range("N19").Value = 0.001523229246001523
$ws.Range("O19").Value = 0.09215536938309216
$ws.Range("S19").Value = 0.07996953541507996
